# Update 'want to go' counts (column F) across sheets per upstream data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 168
$ws.Range("F6").Value = 2749
$ws.Range("F9").Value = 7447
$ws.Range("F11").Value = 7635
$ws.Range("F13").Value = 34
$ws.Range("F15").Value = 6157
$ws.Range("F16").Value = 3254
$ws.Range("F19").Value = 10
$ws.Range("F24").Value = 284
$ws.Range("F25").Value = 285
$ws.Range("F26").Value = 3621
$ws.Range("F28").Value = 340
$ws.Range("F31").Value = 1086
$ws.Range("F34").Value = 2609
$ws.Range("F35").Value = 1458
$ws.Range("F36").Value = 11
$ws.Range("F37").Value = 15
$ws.Range("F39").Value = 3255
$ws.Range("F40").Value = 157
$ws.Range("F44").Value = 480
$ws.Range("F45").Value = 1273
$ws.Range("F47").Value = 523
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 234
$ws.Range("F9").Value = 401
$ws.Range("F10").Value = 30
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 168
$ws.Range("F10").Value = 234
$ws.Range("F12").Value = 7447
$ws.Range("F13").Value = 7635
$ws.Range("F15").Value = 6157
$ws.Range("F16").Value = 3254
$ws.Range("F19").Value = 10
$ws.Range("F23").Value = 284
$ws.Range("F26").Value = 285
$ws.Range("F27").Value = 3621
$ws.Range("F30").Value = 340
$ws.Range("F35").Value = 2609
$ws.Range("F36").Value = 1458
$ws.Range("F37").Value = 11
$ws.Range("F38").Value = 15
$ws.Range("F40").Value = 3255
$ws.Range("F41").Value = 157
$ws.Range("F46").Value = 480
$ws.Range("F47").Value = 1273
$ws.Range("F49").Value = 523
